$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# --- Crime statistics table updates ---
$ws.Range("N14").Value = -100
$ws.Range("N14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 4
$ws.Range("K16").Value = 33.333333333333
$ws.Range("L16").Value = 33.333333333333
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = -50
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 37.5
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = 81.818181818181
$ws.Range("L17").Value = 122.222222222222
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = 185.714285714286
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("I18").Value = 11
$ws.Range("K18").Value = 266.666666666667
$ws.Range("L18").Value = 10
$ws.Range("M18").Value = -31.25
$ws.Range("N18").Value = -74.418604651162
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -62.5
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -39.130434782608
$ws.Range("I19").Value = 29
$ws.Range("J19").Value = 36
$ws.Range("K19").Value = -19.444444444444
$ws.Range("L19").Value = -30.952380952381
$ws.Range("M19").Value = 11.538461538461
$ws.Range("N19").Value = 38.095238095238
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 100
$ws.Range("N20").Value = -95.327102803738
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 40
$ws.Range("G21").Value = 37
$ws.Range("H21").Value = 8.108108108108
$ws.Range("I21").Value = 69
$ws.Range("J21").Value = 62
$ws.Range("K21").Value = 11.290322580645
$ws.Range("L21").Value = -6.756756756756
$ws.Range("M21").Value = 15
$ws.Range("N21").Value = -63.101604278074
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 150
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 29
$ws.Range("H24").Value = 10.344827586206
$ws.Range("I24").Value = 70
$ws.Range("J24").Value = 59
$ws.Range("K24").Value = 18.644067796610
$ws.Range("L24").Value = -11.392405063291
$ws.Range("M24").Value = -10.256410256410
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 6.666666666666
$ws.Range("I25").Value = 33
$ws.Range("J25").Value = 25
$ws.Range("K25").Value = 32
$ws.Range("L25").Value = -21.428571428571
$ws.Range("C26").Value = 7
$ws.Range("E26").Value = 75
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = 108.333333333333
$ws.Range("I26").Value = 37
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = 94.736842105263
$ws.Range("L26").Value = 19.354838709677
$ws.Range("M26").Value = 19.354838709677
$ws.Range("G27").Value = 2
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = 300
$ws.Range("L28").Value = 100
$ws.Range("N29").Value = -100
$ws.Range("N29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N30").Value = -100
$ws.Range("N30").NumberFormat = '#,##0.0;"-"#,##0.0'
